# Natmi following Dr Hou advice
# Ligand/receptor-expressing cell counts increased from 1 to 3 for each
# Sending/Target cluster pair, with all dependent average/total expression,
# specificity and edge-weight statistics recomputed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value2 = 3
$ws.Range("G2").Value2 = 1.191131666666666
$ws.Range("H2").Value2 = 3.573395
$ws.Range("I2").Value2 = 0.02720036629735778
$ws.Range("J2").Value2 = 0.02720036629735778
$ws.Range("K2").Value2 = 3
$ws.Range("M2").Value2 = 4.847498666666667
$ws.Range("N2").Value2 = 14.542496
$ws.Range("O2").Value2 = 0.03400671694637637
$ws.Range("P2").Value2 = 0.03400671694637637
$ws.Range("Q2").Value2 = 5.77400916599111
$ws.Range("R2").Value2 = 51.96608249391999
$ws.Range("S2").Value2 = 0.0009249951575120012
$ws.Range("T2").Value2 = 0.0009249951575120012

# Row 3
$ws.Range("E3").Value2 = 3
$ws.Range("G3").Value2 = 1.191131666666666
$ws.Range("H3").Value2 = 3.573395
$ws.Range("I3").Value2 = 0.02720036629735778
$ws.Range("J3").Value2 = 0.02720036629735778
$ws.Range("K3").Value2 = 3
$ws.Range("M3").Value2 = 1.627093333333334
$ws.Range("N3").Value2 = 4.88128
$ws.Range("O3").Value2 = 0.01141456784970118
$ws.Range("P3").Value2 = 0.01141456784970118
$ws.Range("Q3").Value2 = 1.938082393955555
$ws.Range("R3").Value2 = 17.4427415456
$ws.Range("S3").Value2 = 0.0003104804266379157
$ws.Range("T3").Value2 = 0.0003104804266379157

# Row 4
$ws.Range("E4").Value2 = 3
$ws.Range("G4").Value2 = 1.191131666666666
$ws.Range("H4").Value2 = 3.573395
$ws.Range("I4").Value2 = 0.02720036629735778
$ws.Range("J4").Value2 = 0.02720036629735778
$ws.Range("K4").Value2 = 3
$ws.Range("M4").Value2 = 136.0707373333333
$ws.Range("N4").Value2 = 408.212212
$ws.Range("O4").Value2 = 0.9545787152039225
$ws.Range("P4").Value2 = 0.9545787152039225
$ws.Range("Q4").Value2 = 162.0781641444156
$ws.Range("R4").Value2 = 1458.70347729974
$ws.Range("S4").Value2 = 0.02596489071320786
$ws.Range("T4").Value2 = 0.02596489071320786

# Row 5
$ws.Range("E5").Value2 = 3
$ws.Range("G5").Value2 = 34.415161
$ws.Range("H5").Value2 = 103.245483
$ws.Range("I5").Value2 = 0.7858954736735307
$ws.Range("J5").Value2 = 0.7858954736735306
$ws.Range("K5").Value2 = 3
$ws.Range("M5").Value2 = 4.847498666666667
$ws.Range("N5").Value2 = 14.542496
$ws.Range("O5").Value2 = 0.03400671694637637
$ws.Range("P5").Value2 = 0.03400671694637637
$ws.Range("Q5").Value2 = 166.8274470606187
$ws.Range("R5").Value2 = 1501.447023545568
$ws.Range("S5").Value2 = 0.02672572492265414
$ws.Range("T5").Value2 = 0.02672572492265414

# Row 6
$ws.Range("E6").Value2 = 3
$ws.Range("G6").Value2 = 34.415161
$ws.Range("H6").Value2 = 103.245483
$ws.Range("I6").Value2 = 0.7858954736735307
$ws.Range("J6").Value2 = 0.7858954736735306
$ws.Range("K6").Value2 = 3
$ws.Range("M6").Value2 = 1.627093333333334
$ws.Range("N6").Value2 = 4.88128
$ws.Range("O6").Value2 = 0.01141456784970118
$ws.Range("P6").Value2 = 0.01141456784970118
$ws.Range("Q6").Value2 = 55.99667902869335
$ws.Range("R6").Value2 = 503.9701112582401
$ws.Range("S6").Value2 = 0.008970657207019565
$ws.Range("T6").Value2 = 0.008970657207019565

# Row 7
$ws.Range("E7").Value2 = 3
$ws.Range("G7").Value2 = 34.415161
$ws.Range("H7").Value2 = 103.245483
$ws.Range("I7").Value2 = 0.7858954736735307
$ws.Range("J7").Value2 = 0.7858954736735306
$ws.Range("K7").Value2 = 3
$ws.Range("M7").Value2 = 136.0707373333333
$ws.Range("N7").Value2 = 408.212212
$ws.Range("O7").Value2 = 0.9545787152039225
$ws.Range("P7").Value2 = 0.9545787152039225
$ws.Range("Q7").Value2 = 4682.896332715378
$ws.Range("R7").Value2 = 42146.0669944384
$ws.Range("S7").Value2 = 0.750199091543857
$ws.Range("T7").Value2 = 0.7501990915438569

# Row 8
$ws.Range("E8").Value2 = 3
$ws.Range("G8").Value2 = 8.184723
$ws.Range("H8").Value2 = 24.554169
$ws.Range("I8").Value2 = 0.1869041600291116
$ws.Range("J8").Value2 = 0.1869041600291116
$ws.Range("K8").Value2 = 3
$ws.Range("M8").Value2 = 4.847498666666667
$ws.Range("N8").Value2 = 14.542496
$ws.Range("O8").Value2 = 0.03400671694637637
$ws.Range("P8").Value2 = 0.03400671694637637
$ws.Range("Q8").Value2 = 39.675433829536
$ws.Range("R8").Value2 = 357.078904465824
$ws.Range("S8").Value2 = 0.00635599686621023
$ws.Range("T8").Value2 = 0.006355996866210229

# Row 9
$ws.Range("E9").Value2 = 3
$ws.Range("G9").Value2 = 8.184723
$ws.Range("H9").Value2 = 24.554169
$ws.Range("I9").Value2 = 0.1869041600291116
$ws.Range("J9").Value2 = 0.1869041600291116
$ws.Range("K9").Value2 = 3
$ws.Range("M9").Value2 = 1.627093333333334
$ws.Range("N9").Value2 = 4.88128
$ws.Range("O9").Value2 = 0.01141456784970118
$ws.Range("P9").Value2 = 0.01141456784970118
$ws.Range("Q9").Value2 = 13.31730822848
$ws.Range("R9").Value2 = 119.85577405632
$ws.Range("S9").Value2 = 0.002133430216043702
$ws.Range("T9").Value2 = 0.002133430216043702

# Row 10
$ws.Range("E10").Value2 = 3
$ws.Range("G10").Value2 = 8.184723
$ws.Range("H10").Value2 = 24.554169
$ws.Range("I10").Value2 = 0.1869041600291116
$ws.Range("J10").Value2 = 0.1869041600291116
$ws.Range("K10").Value2 = 3
$ws.Range("M10").Value2 = 136.0707373333333
$ws.Range("N10").Value2 = 408.212212
$ws.Range("O10").Value2 = 0.9545787152039225
$ws.Range("P10").Value2 = 0.9545787152039225
$ws.Range("Q10").Value2 = 1113.701293479092
$ws.Range("R10").Value2 = 10023.31164131183
$ws.Range("S10").Value2 = 0.1784147329468577
$ws.Range("T10").Value2 = 0.1784147329468577
